$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actual")

$rowRange = $ws.Range("A12:H12")
# Force text interpretation so numeric/date-looking values ("10541212",
# "2024-02-07") are stored as literal text instead of being auto-converted
# to a number/date, matching the rest of the sheet's columns.
$rowRange.NumberFormat = "@"

$ws.Cells.Item(12, 1).Value = "10541212"
$ws.Cells.Item(12, 2).Value = "prueba Proveedor 2"
$ws.Cells.Item(12, 3).Value = "Tercero"
$ws.Cells.Item(12, 4).Value = "2024-02-07"
$ws.Cells.Item(12, 5).Value = "16:57:55"
$ws.Cells.Item(12, 6).Value = "16:38:14"
$ws.Cells.Item(12, 7).Value = "16:57:55"
$ws.Cells.Item(12, 8).Value = "Salida PM"

# Reset the visual style back to Normal now that the text values are
# committed, so the new row doesn't carry a stray "@" number-format style
# that the rest of the sheet's rows don't have.
$rowRange.Style = "Normal"
